$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fecha" (date) column header in G1 - plain/default style, like the
# rest of the newly-added column data.
$ws.Range("G1").Value = "Fecha"

# Build the date-stamped style on G2 first (number format + italic green
# Arial font), then stamp the same style across G3:G11 via copy/paste so
# every data row shares one cell style (matches the single new cellXfs
# entry used by the whole column).
$g2 = $ws.Range("G2")
$g2.Value = 42528.5654282407
$g2.NumberFormat = "DD/MM/YYYY\ HH:MM:SS"
$g2.Font.Name = "Arial"
$g2.Font.Size = 10
$g2.Font.Italic = $true
$g2.Font.Color = 43520

$rest = $ws.Range("G3:G11")
$rest.Value = 42528.5654282407
$g2.Copy()
$rest.PasteSpecial(-4122)

# Column widths — nudge to the post-edit widths from the authored file.
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 10.5
$ws.Columns.Item(3).ColumnWidth = 18.3333333333333
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 13
$ws.Columns.Item(6).ColumnWidth = 10.5
$ws.Columns.Item(7).ColumnWidth = 18.6666666666667

# Selection moves to G21, matching the post-edit view state.
$ws.Range("G21").Select() | Out-Null
